# Remove the 2000年-2009年 rows (rows 2-11), shifting the 2010年-2015年
# rows (previously rows 12-17) up so they become the new rows 2-7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:D11").EntireRow.Delete()
